# Apply table edits to Sheet1: update wording/citations, add three new
# data rows (Cramwinckel et al. 2018 / Evans et al. 2018 / Tierney et al. 2017),
# left-align the numeric "value (error)" cells that were introduced, and
# refresh the selection to span the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Full target contents of the table (text cells as strings, numeric
#     cells as numbers) -----------------------------------------------
$data = @(
  @("Source",                      "Time",         "Type", "Gradient (average across both hemispheres)", "uncertainty", "type of gradient",       "model",                 "proxy system",                                 "comment"),
  @("Bij et al. 2009",             "early Eocene", "SST",  "7 (9, 5)",                                    "-",           "equator - polar circle", "2nd order polynomial", "TEX86, Uk37",                                  "excluding outlier d18O datum"),
  @("Keating-Bitoni et al. 2011",  "EECO",         "SST",  "13 (11, 14)",                                 "-",           "equator - polar circle", "2nd order polynomial", "TEX86, MBT, clumped, Mg/Ca, d18O",            $null),
  @("Cramwinckel et al. 2018",     "EECO",         "SST",  21,                                            1,             "equator - deep water",   "none",                  "TEX86, clumped, Mg/Ca, d18O, deepwater d18O", $null),
  @("Evans et al. 2018",           "early Eocene", "SST",  20,                                            3,             "equator - deep water",   "none",                  "clumped, deep-water Mg/Ca",                   "mainly non-EECO data"),
  @("Tierney et al. 2017",         "early Eocene", "SST",  12,                                            "-",           "equator - polar circle", "Gaussian function",    "TEX86",                                        $null)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# --- Newly introduced numeric cells (columns D and E on the new rows,
#     plus D6) are left aligned, matching the rest of the text columns --
$xlLeft = -4131   # Microsoft.Office.Interop.Excel.XlHAlign.xlHAlignLeft
foreach ($addr in @("D4", "E4", "D5", "E5", "D6")) {
    $ws.Range($addr).HorizontalAlignment = $xlLeft
}

# --- Column widths adjusted slightly (table grew wider/taller) --------
$ws.Columns.Item(1).ColumnWidth = 21.5
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 11.5
$ws.Columns.Item(4).ColumnWidth = 45.666666666666664
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666
$ws.Columns.Item(6).ColumnWidth = 20.666666666666668
$ws.Columns.Item(7).ColumnWidth = 19
$ws.Columns.Item(8).ColumnWidth = 29.833333333333332
$ws.Columns.Item(9).ColumnWidth = 24.5

# --- Selection now spans the full table ------------------------------
$ws.Range("A1:I6").Select()
